$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing Excel to keep/store it as
# plain text (matches the workbooks existing text-cell convention),
# then restore the default "Normal" style so no stray formatting is
# left behind on cells whose printed text happens to look numeric
# (e.g. "58.20", "0.999") -- otherwise Excel auto-converts them to
# numbers and trailing zeros / formatting would be lost.
function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "70.844.55"
$ws.Range("E2").Value = "  +7.45%  "

# Row 3
$ws.Range("D3").Value = "3.634.26"
$ws.Range("E3").Value = "  +7.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
Set-TextValue "D5" "594.75"
$ws.Range("E5").Value = "  +5.45%  "

# Row 6
Set-TextValue "D6" "192.31"
$ws.Range("E6").Value = "  +8.95%  "

# Row 7
Set-TextValue "D7" "0.656"
$ws.Range("E7").Value = "  +4.12%  "

# Row 8
$ws.Range("D8").Value = "3.622.64"
$ws.Range("E8").Value = "  +7.14%  "

# Row 9
Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("E10").Value = "  +3.06%  "

# Row 11
$ws.Range("E11").Value = "  +5.15%  "

# Row 12
Set-TextValue "D12" "58.20"
$ws.Range("E12").Value = "  +7.92%  "

# Row 13
$ws.Range("E13").Value = "  +6.04%  "

# Row 15
$ws.Range("D15").Value = "4.212.98"
$ws.Range("E15").Value = "  +7.24%  "

# Row 16
Set-TextValue "D16" "19.46"
$ws.Range("E16").Value = "  +6.88%  "

# Row 17
$ws.Range("D17").Value = "3.625.93"
$ws.Range("E17").Value = "  +7.20%  "

# Row 18
$ws.Range("D18").Value = "70.660.75"
$ws.Range("E18").Value = "  +7.50%  "

# Row 19
Set-TextValue "D19" "12.65"
$ws.Range("E19").Value = "  +6.20%  "

# Row 20
$ws.Range("E20").Value = "  +0.89%  "

# Row 21
Set-TextValue "D21" "1.06"
$ws.Range("E21").Value = "  +5.90%  "

# Row 22
Set-TextValue "D22" "495.45"
$ws.Range("E22").Value = "  +6.89%  "

# Row 23
$ws.Range("E23").Value = "  +10.91%  "

# Row 24
Set-TextValue "D24" "17.11"
$ws.Range("E24").Value = "  +15.30%  "

# Row 25
Set-TextValue "D25" "4.50"
$ws.Range("E25").Value = "  +9.48%  "

# Row 26
Set-TextValue "D26" "91.27"
$ws.Range("E26").Value = "  +2.13%  "

# Row 27
$ws.Range("E27").Value = "  +7.16%  "

# Row 28
Set-TextValue "D28" "11.31"
$ws.Range("E28").Value = "  +6.10%  "

# Row 29
Set-TextValue "D29" "9.49"
$ws.Range("E29").Value = "  +8.81%  "

# Row 30
Set-TextValue "D30" "32.42"
$ws.Range("E30").Value = "  +4.17%  "

# Row 31
Set-TextValue "D31" "7.65"
$ws.Range("E31").Value = "  +15.55%  "

# Row 32
Set-TextValue "D32" "12.28"
$ws.Range("E32").Value = "  +6.89%  "

# Row 33
Set-TextValue "D33" "619.16"
$ws.Range("E33").Value = "  +6.72%  "

# Row 34
$ws.Range("E34").Value = "  +8.76%  "

# Row 35
Set-TextValue "D35" "65.30"
$ws.Range("E35").Value = "  +4.59%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0835"
$ws.Range("E36").Value = "  +11.29%  "

# Row 37
Set-TextValue "D37" "0.409"
$ws.Range("E37").Value = "  +7.88%  "

# Row 38
$ws.Range("E38").Value = "  +3.89%  "

# Row 39
Set-TextValue "D39" "38.23"
$ws.Range("E39").Value = "  +6.04%  "

# Row 40
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
Set-TextValue "D41" "3.66"
$ws.Range("E41").Value = "  +1.97%  "

# Row 42
$ws.Range("D42").Value = "3.342.47"
$ws.Range("E42").Value = "  +7.71%  "

# Row 43
Set-TextValue "D43" "3.09"
$ws.Range("E43").Value = "  +8.25%  "

# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D44" "2.72"
$ws.Range("E44").Value = "  +11.09%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0449"
$ws.Range("E45").Value = "  +7.38%  "

# Row 46
Set-TextValue "D46" "3.37"
$ws.Range("E46").Value = "  +5.35%  "

# Row 47
Set-TextValue "D47" "0.138"
$ws.Range("E47").Value = "  +3.27%  "

# Row 48
Set-TextValue "D48" "9.25"
$ws.Range("E48").Value = "  +8.89%  "

# Row 49
$ws.Range("E49").Value = "  +6.94%  "

# Row 50
Set-TextValue "D50" "3.36"
$ws.Range("E50").Value = "  +5.96%  "

# Row 51
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D51" "0.999"
$ws.Range("E51").Value = "  -0.05%  "
